$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90: politeness_score (column B) changes from text "3" to numeric 3.
# All other cells in row 90 remain unchanged.
$ws.Cells.Item(90, 2).Value = 3

# New row 91: new annotation entry appended below row 90.
$ws.Cells.Item(91, 1).Value = "Ruilin"

# Column B (politeness_score) is stored as text "3" here (not a number),
# so force text entry and strip any auto-applied number formatting.
$ws.Cells.Item(91, 2).Value = "'3"
$ws.Cells.Item(91, 2).ClearFormats()

$ws.Cells.Item(91, 3).Value = "无"
$ws.Cells.Item(91, 4).Value = "DIS"
$ws.Cells.Item(91, 5).Value = "OTH"
$ws.Cells.Item(91, 6).Value = "216e3c96-70ff-4d1d-bc9b-ae161e0068a3"
$ws.Cells.Item(91, 7).Value = "BJInEZsTb_annotated.xlsx"
$ws.Cells.Item(91, 8).Value = "While this is true, we do not believe is necessarily constitutes a disadvantage of our networks, especially when considering ease of training and reproducibility."

Write-Host "done"
